# Optuna Attempt (go back with original)
# Update forecast values on "Forecast Comparison" and recompute Summary stats.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison" sheet ---
# Row 2
$wsForecast.Range("D2").Value = 4
$wsForecast.Range("H2").Value = 1.5
$wsForecast.Range("J2").Value = "Normal"
$wsForecast.Range("L2").Value = 1.13

# Row 3
$wsForecast.Range("D3").Value = 4
$wsForecast.Range("H3").Value = 0.4
$wsForecast.Range("L3").Value = 1.02

# Row 4
$wsForecast.Range("D4").Value = 5
$wsForecast.Range("L4").Value = 0.9399999999999999

# Row 5
$wsForecast.Range("D5").Value = 6
$wsForecast.Range("L5").Value = 1.05

# Row 6
$wsForecast.Range("D6").Value = 6
$wsForecast.Range("L6").Value = 0.96

# Row 7
$wsForecast.Range("D7").Value = 6
$wsForecast.Range("L7").Value = 1.11

# Row 8
$wsForecast.Range("D8").Value = 7
$wsForecast.Range("L8").Value = 1.04

# Row 9
$wsForecast.Range("D9").Value = 7
$wsForecast.Range("L9").Value = 0.96

# Row 10
$wsForecast.Range("D10").Value = 6
$wsForecast.Range("L10").Value = 1.04

# Row 11
$wsForecast.Range("D11").Value = 7
$wsForecast.Range("L11").Value = 0.89

# Row 12
$wsForecast.Range("D12").Value = 7
$wsForecast.Range("L12").Value = 1.03

# Row 13
$wsForecast.Range("D13").Value = 7
$wsForecast.Range("L13").Value = 0.97

# Row 14
$wsForecast.Range("D14").Value = 7
$wsForecast.Range("L14").Value = 1.05

# Row 15
$wsForecast.Range("D15").Value = 6
$wsForecast.Range("L15").Value = 0.97

# Row 16
$wsForecast.Range("D16").Value = 7
$wsForecast.Range("L16").Value = 1.15

# Row 17
$wsForecast.Range("D17").Value = 6
$wsForecast.Range("L17").Value = 0.82

# --- "Summary" sheet ---
# These cells hold numeric-looking labels stored as text, so briefly force
# the Text number format before assigning (prevents Excel's automatic
# numeric coercion of a digit-only string), then restore the plain "Normal"
# style so no formatting footprint is left behind -- only the text value
# actually changes, matching the source workbook's untouched styling.
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "104"
$wsSummary.Range("B9").Style = "Normal"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "48"
$wsSummary.Range("B10").Style = "Normal"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "21"
$wsSummary.Range("B11").Style = "Normal"

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "7"
$wsSummary.Range("B12").Style = "Normal"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "4"
$wsSummary.Range("B14").Style = "Normal"
